$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 currently holds the plain number 20240716 (a "yyyymmdd" style label).
# Convert it into a real Excel date serial (2024-07-16) and apply a date
# number format to it (reuses the existing built-in short-date format,
# numFmtId 14).
$ws.Range("A1").Value = 45489
$ws.Range("A1").NumberFormat = "m/d/yy"

# Move the active selection to A2.
$ws.Range("A2").Select()

$wb.Save()
